# contratos-6-2016.xlsx — "fix: fixed formatting when scrapping floating point numbers"
#
# Two "Razon social" text cells used a comma as a separator between co-contractors;
# those commas are normalized to periods (plain text edit, no numeric parsing involved).
#
# The "Importe" column (H) holds amounts that were scraped using Argentine/es-AR
# formatting ("." thousands separator, "," decimal separator), e.g. "1.945,36". The fix
# re-writes them using plain decimal-point notation ("1945.36") while keeping them as
# literal TEXT cells (they were text before, and must stay text): the numeric value is
# prefixed with an apostrophe so Excel does not auto-convert it into a real number, and
# the range style is reset back to Normal afterwards so no left-over Text number-format
# is applied to the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E79").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E154").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

$importe = @{
    2 = "840.00"
    3 = "800.00"
    4 = "261993.00"
    5 = "300000.00"
    6 = "1945.36"
    7 = "559.00"
    8 = "60.00"
    9 = "11889.00"
    10 = "6050.00"
    11 = "550.00"
    12 = "5150.00"
    13 = "850.00"
    14 = "333545.80"
    15 = "144529.04"
    16 = "53736.50"
    17 = "49293.76"
    18 = "3901.95"
    19 = "34262.00"
    20 = "11624.60"
    21 = "21951.47"
    22 = "21181.21"
    23 = "2520.00"
    24 = "357.05"
    25 = "431.40"
    26 = "1205.00"
    27 = "3576.99"
    28 = "115990.06"
    29 = "940.00"
    30 = "18356.06"
    31 = "1194.80"
    32 = "5471.50"
    33 = "163.36"
    34 = "240.00"
    35 = "29.13"
    36 = "77.40"
    37 = "558.00"
    38 = "13316.60"
    39 = "46161.18"
    40 = "1735.60"
    41 = "7943.15"
    42 = "6746.00"
    43 = "3450.00"
    44 = "9302.97"
    45 = "10883.53"
    46 = "298316.00"
    47 = "2668.00"
    48 = "558.69"
    49 = "3857.65"
    50 = "37.20"
    51 = "207464.11"
    52 = "148.80"
    53 = "96727.66"
    54 = "84700.00"
    55 = "700.00"
    56 = "1704.30"
    57 = "34185.70"
    58 = "2916.00"
    59 = "768.00"
    60 = "1505.44"
    61 = "9550.00"
    62 = "2800.00"
    63 = "14025.00"
    64 = "38700.00"
    65 = "1986.95"
    66 = "10504.97"
    67 = "510.00"
    68 = "1220.00"
    69 = "102000.00"
    70 = "56500.00"
    71 = "500.02"
    72 = "310.00"
    73 = "4400.00"
    74 = "9800.00"
    75 = "8600.00"
    76 = "104.32"
    77 = "13356.00"
    78 = "2540.00"
    79 = "5940.00"
    80 = "2680.00"
    81 = "2380.00"
    82 = "1400.00"
    83 = "9423.00"
    84 = "233898.00"
    85 = "365.00"
    86 = "415230.93"
    87 = "3157.00"
    88 = "53261.49"
    89 = "4500.00"
    90 = "3600.00"
    91 = "29.40"
    92 = "48.68"
    93 = "31467.15"
    94 = "8270.00"
    95 = "146.00"
    96 = "13047.00"
    97 = "22230.00"
    98 = "14700.75"
    99 = "74.70"
    100 = "3597.00"
    101 = "70.00"
    102 = "2582.00"
    103 = "20663.11"
    104 = "5175.00"
    105 = "987.20"
    106 = "1294.36"
    107 = "1165.00"
    108 = "555.70"
    109 = "4550.00"
    110 = "1328.60"
    111 = "1516.00"
    112 = "207717.50"
    113 = "84400.00"
    114 = "1778.00"
    115 = "6900.00"
    116 = "3079.03"
    117 = "5850.00"
    118 = "3400.00"
    119 = "41288.00"
    120 = "7000.00"
    121 = "4000.00"
    122 = "6050.00"
    123 = "354.32"
    124 = "2110.00"
    125 = "6478.62"
    126 = "1637.40"
    127 = "625.00"
    128 = "450.52"
    129 = "301500.00"
    130 = "2068.56"
    131 = "5125.00"
    132 = "7500.00"
    133 = "200.00"
    134 = "2880.00"
    135 = "2178.00"
    136 = "6200.00"
    137 = "5850.00"
    138 = "4830.00"
    139 = "4900.00"
    140 = "109.64"
    141 = "67.10"
    142 = "26000.00"
    143 = "8190.00"
    144 = "545.69"
    145 = "266.00"
    146 = "733.20"
    147 = "2020.00"
    148 = "2623.08"
    149 = "22400.00"
    150 = "212.80"
    151 = "1100.00"
    152 = "1521.00"
    153 = "5220.17"
    154 = "22600.00"
    155 = "8200.00"
    156 = "3372.00"
    157 = "632.00"
    158 = "319.80"
    159 = "33540.00"
    160 = "6492.00"
    161 = "3400.00"
    162 = "4650.00"
    163 = "461.58"
    164 = "5326.27"
    165 = "14369.95"
    166 = "25000.00"
    167 = "25000.00"
    168 = "25000.00"
    169 = "25000.00"
    170 = "25000.00"
    171 = "25000.00"
    172 = "7000.00"
    173 = "169400.00"
    174 = "4219992.55"
    175 = "1500.00"
    176 = "141583.00"
    177 = "105000.00"
    178 = "105000.00"
    179 = "105000.00"
    180 = "155000.00"
    181 = "105000.00"
    182 = "175000.00"
    183 = "200000.00"
    184 = "245000.00"
    185 = "105000.00"
    186 = "105000.00"
    187 = "105000.00"
    188 = "105000.00"
    189 = "130000.00"
    190 = "200000.00"
    191 = "297000.00"
    192 = "200000.00"
    193 = "105000.00"
    194 = "130000.00"
    195 = "105000.00"
    196 = "105000.00"
    197 = "130000.00"
    198 = "210103.61"
    199 = "105000.00"
    200 = "8223.00"
    201 = "11796.00"
    202 = "14700.00"
    203 = "13506.14"
    204 = "12864.22"
    205 = "54448.80"
    206 = "19337.50"
    207 = "17500.00"
    208 = "725.00"
}

foreach ($row in $importe.Keys) {
    $ws.Cells.Item([int]$row, 8).Value = "'" + $importe[$row]
}

# Drop the quote-prefix formatting the loop above applied so the cells end up with the
# same (default) style they started with — only their text content changed.
$ws.Range("H2:H208").Style = "Normal"

